$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 96
$ws.Cells.Item(96, 1).Value = 9
$ws.Cells.Item(96, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(96, 3).Value = "Metropolitana"
$ws.Cells.Item(96, 4).Value = 44890
$ws.Cells.Item(96, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(96, 5).Value = 13
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100103
$ws.Cells.Item(96, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(96, 9).Value = 100103003
$ws.Cells.Item(96, 10).Value = "Damasco"
$ws.Cells.Item(96, 11).Value = "Castle Brite"
$ws.Cells.Item(96, 12).Value = "Primera"
$ws.Cells.Item(96, 13).Value = 250
$ws.Cells.Item(96, 14).Value = 16000
$ws.Cells.Item(96, 15).Value = 16000
$ws.Cells.Item(96, 16).Value = 16000
$ws.Cells.Item(96, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(96, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(96, 19).Value = 1600
$ws.Cells.Item(96, 20).Value = 10

# Row 97
$ws.Cells.Item(97, 1).Value = 9
$ws.Cells.Item(97, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(97, 3).Value = "Metropolitana"
$ws.Cells.Item(97, 4).Value = 44890
$ws.Cells.Item(97, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(97, 5).Value = 13
$ws.Cells.Item(97, 6).Value = "Fruta"
$ws.Cells.Item(97, 7).Value = 100103
$ws.Cells.Item(97, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(97, 9).Value = 100103003
$ws.Cells.Item(97, 10).Value = "Damasco"
$ws.Cells.Item(97, 11).Value = "Castle Brite"
$ws.Cells.Item(97, 12).Value = "Primera"
$ws.Cells.Item(97, 13).Value = 200
$ws.Cells.Item(97, 14).Value = 24000
$ws.Cells.Item(97, 15).Value = 24000
$ws.Cells.Item(97, 16).Value = 24000
$ws.Cells.Item(97, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(97, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(97, 19).Value = 1600
$ws.Cells.Item(97, 20).Value = 15
